# AI4SMM flagship proposal template -- budget paragraph update.
#
# Original sentence (single run):
#   "within the RPA an amount of k€ 1,190 available. The tariff for a PhD
#    position is k€ ***/ 4 years, the tariff for a PD position is k€ ***/
#    1 year. These tariffs can be used for cash contributions from
#    external funds too. In kind contributions from partners should be
#    mentioned separately) "
#
# New sentence:
#   "within the RPA an amount of k€ 1,190 available. The tariff for a PhD
#    position is k€ 360/ 4 years, with k€ 60/20 additional for
#    experimental/computational work. These tariffs can be used for cash
#    contributions from external funds too. In kind contributions from
#    partners should be mentioned separately) "
#
# The author typed the new numbers/clause in as separate fragments, so
# the final document has that clause split across several runs (all
# sharing the same Times New Roman / 0070C0 / 11pt formatting that the
# surrounding text already used).

$d = $word.ActiveDocument
$d.TrackRevisions = $false

# Step 1: replace the two "***" placeholders and the "PD position" clause
# with the new wording in one shot so the paragraph's text content ends
# up exactly right.
$old = "***/ 4 years, the tariff for a PD position is k€ ***/ 1 year"
$new = "360/ 4 years, with k€ 60/20 additional for experimental/computational work"

$hit = $d.Content
$hit.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Step 2: re-apply (identical) run formatting to each of the new
# fragments in left-to-right order so Word splits them into their own
# runs, matching how they were typed/pasted in.
$cursor = $d.Content
$cursor.Start = $hit.Start
$cursor.End = $hit.Start
$docEnd = $d.Content.End

$fragments = @(
  "360",
  "/ 4 years",
  ", with ",
  "k€",
  " 60/20 additional for experimental/computational work"
)

foreach ($frag in $fragments) {
  $cursor.End = $docEnd
  $cursor.Find.Execute($frag, $true, $false, $false, $false, $false, $true, 1, $false)
  $cursor.Font.Name = "Times New Roman"
  $cursor.Font.NameBi = "Times New Roman"
  $cursor.Font.Color = 12611584
  $cursor.Font.Size = 11
  $cursor.Collapse(0)
}
